$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N29").Value = "2025-10-21"
$ws.Range("Q29").Value = 2.26

$ws.Range("N30").Value = "2025-10-21"
$ws.Range("Q30").Value = 2.28
$ws.Range("R30").Value = 2.26
$ws.Range("S30").Value = 2.27
$ws.Range("T30").Value = 2.28
$ws.Range("U30").Value = 2.29

$ws.Range("N48").Value = "2025-10-20"
$ws.Range("Q48").Value = 3.46
$ws.Range("R48").Value = 3.46
$ws.Range("S48").Value = 3.41
$ws.Range("T48").Value = 3.5
$ws.Range("U48").Value = 3.48

$ws.Range("N49").Value = "2025-10-20"
$ws.Range("Q49").Value = 3.58
$ws.Range("R49").Value = 3.59
$ws.Range("S49").Value = 3.55
$ws.Range("T49").Value = 3.63
$ws.Range("U49").Value = 3.6

$ws.Range("N50").Value = "2025-10-20"
$ws.Range("Q50").Value = 4
$ws.Range("R50").Value = 4.02
$ws.Range("S50").Value = 3.99
$ws.Range("T50").Value = 4.05
$ws.Range("U50").Value = 4.03
